# ProjectGantt.xlsx -- "Weekly" sheet gets a second "Next" list (columns H:K)
# next to the existing one (columns F), three stray items move out of the old
# list, and column headers / selections shift accordingly. "Gantt" sheet only
# has its remembered selection changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly")

# Drop the three tasks that leave the F-column "Next" list
$ws.Range("F7").ClearContents()
$ws.Range("F8").ClearContents()
$ws.Range("F9").ClearContents()

# F1 header relabelled, new H1 header added, F6 task text updated
$ws.Range("F1").Value = 'Week3'
$ws.Range("H1").Value = 'Week4'
$ws.Range("F6").Value = 'Js/Ajax moving'

# New H:K "Next" task list, rows 3-33
$ws.Range("H3").Value = 'Stored Proc Validation'
$ws.Range("I3").Value = 'ME/MP'
$ws.Range("J3").Value = 'Done'
$ws.Range("H4").Value = 'Slight code modifications - hit wall?'
$ws.Range("I4").Value = 'LV'
$ws.Range("J4").Value = 'Done'
$ws.Range("H5").Value = 'Input forms'
$ws.Range("I5").Value = 'LV'
$ws.Range("J5").Value = 'Done'
$ws.Range("H6").Value = 'Test game logic'
$ws.Range("I6").Value = 'All'
$ws.Range("J6").Value = 'Done'
$ws.Range("H7").Value = 'Note all actions done in db - tables, views, index etc Ensure have code'
$ws.Range("I7").Value = 'CC'
$ws.Range("J7").Value = 'Done'
$ws.Range("H8").Value = 'Db system diagram (include views)'
$ws.Range("I8").Value = 'MP'
$ws.Range("J8").Value = 'Done'
$ws.Range("H9").Value = 'Hardcode 10 moves'
$ws.Range("I9").Value = 'LV'
$ws.Range("J9").Value = 'Thu'
$ws.Range("K9").Value = 'maybe drive db in future'
$ws.Range("H10").Value = 'Popup saying success, click next level'
$ws.Range("I10").Value = 'LV'
$ws.Range("J10").Value = 'Thu'
$ws.Range("H11").Value = 'score input in python'
$ws.Range("I11").Value = 'ME'
$ws.Range("J11").Value = 'Thu'
$ws.Range("H12").Value = 'jinga score param,levelid,startx,starty, max levels'
$ws.Range("I12").Value = 'ME'
$ws.Range("J12").Value = 'Thu'
$ws.Range("H13").Value = 'extend template for jinja'
$ws.Range("I13").Value = 'ME'
$ws.Range("J13").Value = 'Thu'
$ws.Range("H14").Value = '100 seconds per level - countdown'
$ws.Range("I14").Value = 'LV'
$ws.Range("J14").Value = 'Thu'
$ws.Range("H15").Value = 'Total score on page'
$ws.Range("I15").Value = 'LV'
$ws.Range("J15").Value = 'Thu'
$ws.Range("H16").Value = 'Congratulations, you''ve finished - no score submit'
$ws.Range("I16").Value = 'LV'
$ws.Range("J16").Value = 'Thu'
$ws.Range("K16").Value = 'possibly be interval from db'
$ws.Range("H17").Value = 'HTML page design'
$ws.Range("I17").Value = 'MP'
$ws.Range("J17").Value = 'Sun'
$ws.Range("H18").Value = 'Score board - db table'
$ws.Range("I18").Value = 'CC'
$ws.Range("J18").Value = 'Thu'
$ws.Range("H19").Value = 'SQL Script for adding score with variables'
$ws.Range("I19").Value = 'CC'
$ws.Range("J19").Value = 'Thu'
$ws.Range("H20").Value = 'add total score to db python - new entry  point (addScore), add method in db'
$ws.Range("I20").Value = 'ME'
$ws.Range("J20").Value = 'Sun'
$ws.Range("H21").Value = 'Finish - - enter name, submit score - js/ajax submit'
$ws.Range("I21").Value = 'LV'
$ws.Range("J21").Value = 'Mon'
$ws.Range("H22").Value = 'SQL to display top 10 scores'
$ws.Range("I22").Value = 'CC'
$ws.Range("J22").Value = 'Thu'
$ws.Range("H23").Value = 'Python generate top 10 scores and add to html - jinga'
$ws.Range("I23").Value = 'MP/CC'
$ws.Range("J23").Value = 'Sun'
$ws.Range("H24").Value = 'Remove diamond'
$ws.Range("I24").Value = 'CC'
$ws.Range("J24").Value = 'Thu'
$ws.Range("H25").Value = 'Level 5'
$ws.Range("I25").Value = 'CC'
$ws.Range("J25").Value = 'Thu'
$ws.Range("H26").Value = 'annimations'
$ws.Range("I26").Value = 'Monday'
$ws.Range("J26").Value = 'Mon'
$ws.Range("H27").Value = 'lightning'
$ws.Range("I27").Value = 'Monday'
$ws.Range("J27").Value = 'Mon'
$ws.Range("H28").Value = 'DB - Spatial Index - fix bug'
$ws.Range("I28").Value = 'CC'
$ws.Range("J28").Value = 'Mon'
$ws.Range("H29").Value = 'Unit Tests - fix Icon'
$ws.Range("I29").Value = 'MP'
$ws.Range("J29").Value = 'Mon'
$ws.Range("H30").Value = 'Code comments - other code'
$ws.Range("I30").Value = 'MP'
$ws.Range("H31").Value = 'E-R diagram - add total score'
$ws.Range("I31").Value = 'CC'
$ws.Range("H32").Value = 'Total Score table - db diagram'
$ws.Range("I32").Value = 'MP'
$ws.Range("H33").Value = 'Developers page/roles'
$ws.Range("I33").Value = 'MP'

# Restore the remembered selections on both sheets
$wsGantt = $wb.Worksheets.Item("Gantt")
$wsGantt.Activate()
$wsGantt.Range("A57").Select()
$ws.Activate()
$ws.Range("K24").Select()
